# Insert a new data row at row 155, pushing existing rows 155-228 down to
# 156-229, and populate the new row with the latest price-report entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 155 and below down by one row.
$ws.Rows.Item(155).Insert()

# Populate the newly inserted row 155 with the new record.
$ws.Cells.Item(155, 1).Value = 3
$ws.Cells.Item(155, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(155, 3).Value = "Coquimbo"
$ws.Cells.Item(155, 4).Value = 44489
$ws.Cells.Item(155, 5).Value = 5
$ws.Cells.Item(155, 6).Value = 100112031
$ws.Cells.Item(155, 7).Value = "Poroto verde"
$ws.Cells.Item(155, 8).Value = "Magnum"
$ws.Cells.Item(155, 9).Value = "Primera"
$ws.Cells.Item(155, 10).Value = 38
$ws.Cells.Item(155, 11).Value = 38000
$ws.Cells.Item(155, 12).Value = 38000
$ws.Cells.Item(155, 13).Value = 38000
$ws.Cells.Item(155, 14).Value = "$/malla 25 kilos"
$ws.Cells.Item(155, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(155, 16).Value = 1520
$ws.Cells.Item(155, 17).Value = 25
$ws.Cells.Item(155, 18).Value = "Hortaliza"
